$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19

$ws.Range("B19").Value = "Process inventory"
$ws.Range("B20").Value = "Unit process inventory"

$ws.Range("E19").Value = "PI"
$ws.Range("E20").Value = "UPI"

$ws.Range("C19").Value = "flows entering and leaving process, for LCI datasets"
$ws.Range("C20").Value = "flows entering and leaving process, normalized, for LCI datasets"

$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 1

$ws.Columns.Item(3).ColumnWidth = 53

$ws.Range("B22").Select()
